# Add data for 2023-06-22
# Updates the CTA violent-crime year-to-date workbook with newly
# reported/reclassified incidents across several years and neighborhoods.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("H2").Value = 57
$ws.Range("D3").Value = 70
$ws.Range("F3").Value = 64
$ws.Range("G3").Value = 56
$ws.Range("J3").Value = 89
$ws.Range("C6").Value = 225
$ws.Range("F6").Value = 244
$ws.Range("G6").Value = 228
$ws.Range("H6").Value = 199
$ws.Range("I6").Value = 267
$ws.Range("C7").Value = 304
$ws.Range("D7").Value = 320
$ws.Range("F7").Value = 347
$ws.Range("G7").Value = 337
$ws.Range("H7").Value = 311
$ws.Range("I7").Value = 416
$ws.Range("J7").Value = 352

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("G4").Value = 3
$ws.Range("G9").Value = 3
$ws.Range("J19").Value = 5
$ws.Range("F25").Value = 6
$ws.Range("H25").Value = 7
$ws.Range("D26").Value = 22
$ws.Range("H26").Value = 28
$ws.Range("J26").Value = 13
$ws.Range("F49").Value = 4
$ws.Range("J59").Value = 2
$ws.Range("C61").Value = 2
$ws.Range("G66").Value = 2
$ws.Range("I67").Value = 8
$ws.Range("G73").Value = 9
$ws.Range("C94").Value = 304
$ws.Range("D94").Value = 320
$ws.Range("F94").Value = 347
$ws.Range("G94").Value = 337
$ws.Range("H94").Value = 311
$ws.Range("I94").Value = 416
$ws.Range("J94").Value = 352

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("G5").Value = 7
$ws.Range("G6").Value = 9

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J4").Value = 4
$ws.Range("J5").Value = 5

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 3

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("H2").Value = 9
$ws.Range("D3").Value = 10
$ws.Range("J3").Value = 3
$ws.Range("D6").Value = 22
$ws.Range("H6").Value = 28
$ws.Range("J6").Value = 13

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 3

$ws = $wb.Worksheets.Item('New City')
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 2

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("F4").Value = 4
$ws.Range("H4").Value = 7
$ws.Range("F5").Value = 6
$ws.Range("H5").Value = 7

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("D3").Value = 1
$ws.Range("D5").Value = 4

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("H4").Value = 7
$ws.Range("H5").Value = 8

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("D3").Value = 2
$ws.Range("D5").Value = 2
